# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de handback-status rows.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 22:57:44"
$wsZhCn.Range("H2").Value = "2016-03-20 22:58:06"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 22:57:48"
$wsDeDe.Range("H2").Value = "2016-03-20 22:58:12"
